$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(46073, 3.901, 0.985),
    @(46073.01041666666, 0.377, 0.414),
    @(46073.02083333334, 18.573, 0),
    @(46073.03125, 20.176, 0),
    @(46073.04166666666, 6.017, 0),
    @(46073.05208333334, 6.662, 0.023),
    @(46073.0625, 22.576, 0),
    @(46073.07291666666, 6.944, 0.5669999999999999),
    @(46073.08333333334, 1.151, 5.161),
    @(46073.09375, 0, 3.272),
    @(46073.10416666666, 0, 1.363),
    @(46073.11458333334, 0, 1.529),
    @(46073.125, 0, 3.005),
    @(46073.13541666666, 0, 18.493),
    @(46073.14583333334, 0.237, 1.013),
    @(46073.15625, 13.681, 0.034),
    @(46073.16666666666, 10.317, 0),
    @(46073.17708333334, 3.049, 0),
    @(46073.1875, 17.883, 0),
    @(46073.19791666666, 38.221, 0),
    @(46073.20833333334, 8.974, 0.005),
    @(46073.21875, 29.043, 0),
    @(46073.22916666666, 42.466, 0),
    @(46073.23958333334, 27.872, 0),
    @(46073.25, 7.35, 9.836),
    @(46073.26041666666, 0.425, 1.064),
    @(46073.27083333334, 0.001, 18.562),
    @(46073.28125, 0, 10.846),
    @(46073.29166666666, 3.012, 2.103),
    @(46073.30208333334, 0, 0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

$wb.Save()